# New weekly TSA data point: insert a fresh row 2 (pushing the existing
# history down by one row) and fill it in with the latest date/value,
# matching the same "% change vs year-ago" formula pattern used by every
# other row in the sheet (compares against the row 12 below it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2; everything below shifts
# down by one (row 2 -> row 3, ..., row 78 -> row 79) and relative formula
# references re-anchor automatically.
$ws.Rows.Item(2).EntireRow.Insert()

# The freshly-inserted row 2 comes back with "General" formatting (it
# copied the mostly-blank row 1 above it). Pull the date format (column A)
# and the numeric format (column H) back from row 3 below, which still
# carries the formatting that used to belong to the old row 2.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new data point.
$ws.Range("A2").Value = 45838
$ws.Range("B2").Value = 81310942
# Column B cells elsewhere in the sheet carry no explicit style even though
# the column default style is 1 (centered) - match that by resetting B2
# back to the plain "Normal" style after the value write.
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Formula = "=(B2/B14-1)*100"
